# Regenerate orders with updated distance/sizes.
# The stimulus naming scheme encodes a viewing "Distance" (D51/D64/D80) and
# a face "Size" (S20/S25/S30) inside many shared strings (condition names,
# left/right filenames, and the standalone Distance/Size lookup values).
# This commit renumbers those codes: D51->D55, D64->D69, D80->D86, S30->S31.
# Apply the substitution to every text cell in the used range in one bulk
# read/modify/write pass (Range.Find/Replace is not available here).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$arr = $used.Value2

$rows = $arr.GetLength(0)
$cols = $arr.GetLength(1)

for ($r = 1; $r -le $rows; $r++) {
  for ($c = 1; $c -le $cols; $c++) {
    $v = $arr[$r, $c]
    if ($v -is [string]) {
      $nv = $v.Replace("D51", "D55").Replace("D64", "D69").Replace("D80", "D86").Replace("S30", "S31")
      if ($nv -ne $v) {
        $arr[$r, $c] = $nv
      }
    }
  }
}

$used.Value2 = $arr
